# Applies the "Updated cryptos list" diff to the active sheet.
# D/E-column values are prefixed with a literal apostrophe so Excel
# stores them as text (matching the source file) instead of auto-
# converting numeric-looking strings (e.g. "600.31", "1.00") into
# numbers, which would drop significant trailing zeros/formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'69.000.03"
$ws.Range("E2").Value = "'  -0.27%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "'3.804.18"
$ws.Range("E3").Value = "'  +1.63%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "'  -0.08%  "

# Row 5: BNB
$ws.Range("D5").Value = "'600.31"
$ws.Range("E5").Value = "'  -0.41%  "

# Row 6: Solana
$ws.Range("D6").Value = "'163.63"
$ws.Range("E6").Value = "'  -2.80%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "'3.804.24"
$ws.Range("E7").Value = "'  +1.71%  "

# Row 8: USDC
$ws.Range("E8").Value = "'  +0.14%  "

# Row 9: XRP
$ws.Range("E9").Value = "'  -0.52%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "'  +1.50%  "

# Row 11: Toncoin
$ws.Range("D11").Value = "'6.30"
$ws.Range("E11").Value = "'  -1.66%  "

# Row 12: Cardano
$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "'  -0.35%  "

# Row 13: Avalanche
$ws.Range("D13").Value = "'37.08"
$ws.Range("E13").Value = "'  -2.72%  "

# Row 14: ShibaInu
$ws.Range("D14").Value = "'0.0000245"
$ws.Range("E14").Value = "'  -1.02%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'4.438.48"
$ws.Range("E15").Value = "'  +1.67%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "'3.779.54"
$ws.Range("E16").Value = "'  +1.05%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "'69.137.69"
$ws.Range("E17").Value = "'  -0.09%  "

# Row 18: Polkadot
$ws.Range("D18").Value = "'7.48"
$ws.Range("E18").Value = "'  +2.25%  "

# Row 19: Uniswap
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.114"
$ws.Range("E19").Value = "'  +0.22%  "

# Row 20: TRON
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'11.51"
$ws.Range("E20").Value = "'  +5.48%  "

# Row 21: Chainlink
$ws.Range("D21").Value = "'17.26"
$ws.Range("E21").Value = "'  +1.16%  "

# Row 22: BitcoinCash
$ws.Range("D22").Value = "'485.12"
$ws.Range("E22").Value = "'  -1.95%  "

# Row 23: Polygon
$ws.Range("D23").Value = "'0.719"
$ws.Range("E23").Value = "'  -0.99%  "

# Row 24: PEPE
$ws.Range("D24").Value = "'0.0000160"
$ws.Range("E24").Value = "'  +6.16%  "

# Row 25: Litecoin
$ws.Range("D25").Value = "'84.59"
$ws.Range("E25").Value = "'  -0.29%  "

# Row 26: Fetch.AI
$ws.Range("E26").Value = "'  -2.93%  "

# Row 27: InternetComputer(DFINITY)
$ws.Range("D27").Value = "'12.20"
$ws.Range("E27").Value = "'  -0.92%  "

# Row 28: RenderToken
$ws.Range("D28").Value = "'10.02"
$ws.Range("E28").Value = "'  -1.59%  "

# Row 29: Dai
$ws.Range("E29").Value = "'  -0.14%  "

# Row 30: PancakeSwap
$ws.Range("D30").Value = "'2.96"
$ws.Range("E30").Value = "'  -1.18%  "

# Row 31: NEARProtocol
$ws.Range("D31").Value = "'8.04"
$ws.Range("E31").Value = "'  -1.16%  "

# Row 32: ImmutableX
$ws.Range("E32").Value = "'  -5.09%  "

# Row 33: WrappedeETH
$ws.Range("D33").Value = "'3.961.40"
$ws.Range("E33").Value = "'  +1.92%  "

# Row 34: EthereumClassic
$ws.Range("D34").Value = "'31.69"
$ws.Range("E34").Value = "'  +0.15%  "

# Row 35: RenzoRestakedETH
$ws.Range("D35").Value = "'3.748.29"
$ws.Range("E35").Value = "'  +1.99%  "

# Row 36: Hedera
$ws.Range("E36").Value = "'  -1.70%  "

# Row 37: Mantle
$ws.Range("E37").Value = "'  +0.64%  "

# Row 38: Kaspa
$ws.Range("D38").Value = "'0.139"
$ws.Range("E38").Value = "'  +4.45%  "

# Row 39: Filecoin
$ws.Range("D39").Value = "'5.86"
$ws.Range("E39").Value = "'  -0.27%  "

# Row 40: FirstDigitalUSD
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "'  +0.06%  "

# Row 41: dogwifhat
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.318"
$ws.Range("E41").Value = "'  -1.97%  "

# Row 42: TheGraph
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'3.03"
$ws.Range("E42").Value = "'  +0.38%  "

# Row 43: Bittensor
$ws.Range("D43").Value = "'438.05"
$ws.Range("E43").Value = "'  +0.88%  "

# Row 45: Stacks
$ws.Range("E45").Value = "'  -0.75%  "

# Row 47: Cosmos
$ws.Range("D47").Value = "'8.37"
$ws.Range("E47").Value = "'  -1.30%  "

# Row 48: Maker
$ws.Range("D48").Value = "'2.824.74"
$ws.Range("E48").Value = "'  +1.52%  "

# Row 49: Monero
$ws.Range("D49").Value = "'141.85"
$ws.Range("E49").Value = "'  +0.83%  "

# Row 50: Arweave
$ws.Range("D50").Value = "'39.20"
$ws.Range("E50").Value = "'  -3.07%  "

# Row 51: VeChain
$ws.Range("D51").Value = "'0.0351"
$ws.Range("E51").Value = "'  -0.50%  "
